$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SanityTC")
$ws.Activate()

# "smx dp buttons added" - three new columns (secondarylanguage, description,
# textbox) were appended after the existing last column (primarylanguage) on
# the SanityTC sheet, for the new Sanity_TC2 (SMx DP) row.

# Header row (row 1)
$ws.Range("AW1").Value = "secondarylanguage"
$ws.Range("AX1").Value = "description"
$ws.Range("AY1").Value = "textbox"

# Data row (row 3) - values for the new SMx DP test case
$ws.Range("AW3").Value = "Spanish"
$ws.Range("AX3").Value = "Please fill the survey and help us improve!"
$ws.Range("AY3").Value = "Enter your Name"

# Match formatting of the neighbouring existing columns (AV) for the new
# header cells, the blank row-2 cells, and the row-3 data cells.
$ws.Range("AV1").Copy()
$ws.Range("AW1:AY1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("AV2").Copy()
$ws.Range("AW2:AY2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("AV3").Copy()
$ws.Range("AW3:AY3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the view/selection to reflect the newly added columns
$ws.Range("AY4").Select()
